# Lab7 Rubric - "Updated the instructions and rubric"
# Reworks the Rubric sheet to a single-part ("Forms"/"Raffle contest form")
# layout, renumbers points, adds a bottom border under the new Syntax-and-
# style score, and relabels the Student Score sheet to match.

$wb = $excel.ActiveWorkbook
$rubric = $wb.Worksheets.Item("Rubric")
$score  = $wb.Worksheets.Item("Student Score")

# ---------------------------------------------------------------------
# Rubric sheet
# ---------------------------------------------------------------------

# Header block
$rubric.Range("A1").Value = "Forms"
$rubric.Range("A2").Value = "Raffle contest form"
$rubric.Range("D2").Value = $null

$rubric.Range("B3").Value = "Points"
$rubric.Range("C3").Value = $null
$rubric.Range("E3").Value = $null
$rubric.Range("F3").Value = $null

# Row 4 - "Set form backbround color" (was the heading row), single column of points
$rubric.Range("A4").Value = "Set form backbround color"
$rubric.Range("B4").Value = 1
$rubric.Range("C4").Value = $null
$rubric.Range("E4").Value = $null
$rubric.Range("F4").Value = $null

$rubric.Range("A5").Value = "Form input elements:"

$rubric.Range("A6").Value = "Name"
$rubric.Range("B6").Value = 4
$rubric.Range("C6").Value = $null
$rubric.Range("E6").Value = $null
$rubric.Range("F6").Value = $null

$rubric.Range("A7").Value = "Mailing address"
$rubric.Range("B7").Value = 4
$rubric.Range("C7").Value = $null
$rubric.Range("E7").Value = $null
$rubric.Range("F7").Value = $null

$rubric.Range("A8").Value = "Phone number"
$rubric.Range("B8").Value = 4
$rubric.Range("C8").Value = $null
$rubric.Range("E8").Value = $null
$rubric.Range("F8").Value = $null

$rubric.Range("A9").Value = "Text area"
$rubric.Range("B9").Value = 4
$rubric.Range("C9").Value = $null
$rubric.Range("E9").Value = $null
$rubric.Range("F9").Value = $null

$rubric.Range("A10").Value = "Two field sets and legends"
$rubric.Range("B10").Value = 5
$rubric.Range("C10").Value = $null
$rubric.Range("E10").Value = $null
$rubric.Range("F10").Value = $null

$rubric.Range("A11").Value = "Labels for all input elements"
$rubric.Range("B11").Value = 5
$rubric.Range("C11").Value = $null
$rubric.Range("E11").Value = $null
$rubric.Range("F11").Value = $null

$rubric.Range("A12").Value = "Form footer"
$rubric.Range("B12").Value = 3
$rubric.Range("B12").Font.Bold = $rubric.Range("B12").Font.Bold
$rubric.Range("C12").Value = $null
$rubric.Range("E12").Value = $null
$rubric.Range("F12").Value = $null

$rubric.Range("A13").Value = "Embedded CSS"

$rubric.Range("A14").Value = "Set background color for fields (a different color)"
$rubric.Range("B14").Value = 3
$rubric.Range("C14").Value = $null
$rubric.Range("E14").Value = $null
$rubric.Range("F14").Value = $null

$rubric.Range("A15").Value = "heading (like <h1>) in header for page"
$rubric.Range("B15").Value = 4
$rubric.Range("C15").Value = $null
$rubric.Range("E15").Value = $null
$rubric.Range("F15").Value = $null

$rubric.Range("A16").Value = "Syntax and style"
$rubric.Range("B16").Value = 3
$rubric.Range("C16").Value = $null
$rubric.Range("E16").Value = $null
$rubric.Range("F16").Value = $null
# New bottom border under the last score in this column
$rubric.Range("B16").Borders.Item(9).LineStyle = 1
$rubric.Range("B16").Borders.Item(9).Weight = 2

$rubric.Range("A17").Value = "Total"
$rubric.Range("C17").Value = $null
$rubric.Range("D17").Font.Italic = $true
$rubric.Range("E17").Value = $null
$rubric.Range("F17").Value = $null

# Second (now unused) mini rubric block - clear out its labels/values but
# keep the cell formatting in place
$rubric.Range("E19").Value = $null

$rubric.Range("A20").Value = $null
$rubric.Range("E20").Value = $null

$rubric.Range("A21").Value = $null
$rubric.Range("B21").Value = $null
$rubric.Range("C21").Value = $null
$rubric.Range("E21").Value = $null

$rubric.Range("A24").Value = $null
$rubric.Range("B24").Value = $null
$rubric.Range("C24").Value = $null

$rubric.Range("A25").Value = $null
$rubric.Range("B25").Value = $null
$rubric.Range("C25").Value = $null

# Column widths: A widened for the longer labels, new narrow B column for points
$rubric.Columns.Item(1).ColumnWidth = 43.166666666666664
$rubric.Columns.Item(2).ColumnWidth = 5.166666666666667

# ---------------------------------------------------------------------
# Student Score sheet - relabel to match the trimmed rubric wording
# ---------------------------------------------------------------------

$score.Range("A1").Value = "We only did part 1 for fall term, 2019"
$score.Range("A3").Value = "Set form backbround color"
$score.Range("A4").Value = "Form input elements:"
$score.Range("A5").Value = "Name"
$score.Range("A6").Value = "Mailing address"
$score.Range("A7").Value = "Phone number"
$score.Range("A8").Value = "Text area"
$score.Range("A9").Value = "Two field sets and legends"
$score.Range("A10").Value = "Labels for all input elements"
$score.Range("A11").Value = "Form footer"
$score.Range("A12").Value = "Embedded CSS"
$score.Range("A13").Value = "Set background color for fields (a different color)"
$score.Range("A14").Value = "heading (like <h1>) in header for page"
$score.Range("A15").Value = "Syntax and style"

# ---------------------------------------------------------------------
# Active sheet / selection - Rubric tab now the one shown, on D9
# ---------------------------------------------------------------------

$rubric.Activate()
$rubric.Range("D9").Select()
